# Applies a weekly refresh of the "Rabanito" dataset: the values in columns
# D (Fecha), I (Calidad), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), O (Origen) and P (Precio $/Kg) are
# redistributed across data rows 2-43 according to a fixed permutation
# (row 36 is unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (values are taken from the
# *original* contents of the source row before any writes happen).
$rowMap = @{
    2  = 13
    3  = 16
    4  = 30
    5  = 15
    6  = 31
    7  = 20
    8  = 21
    9  = 19
    10 = 22
    11 = 41
    12 = 8
    13 = 25
    14 = 26
    15 = 23
    16 = 3
    17 = 12
    18 = 2
    19 = 37
    20 = 43
    21 = 9
    22 = 35
    23 = 24
    24 = 28
    25 = 40
    26 = 42
    27 = 39
    28 = 14
    29 = 4
    30 = 7
    31 = 11
    32 = 6
    33 = 38
    34 = 18
    35 = 27
    36 = 36
    37 = 34
    38 = 17
    39 = 29
    40 = 10
    41 = 32
    42 = 33
    43 = 5
}

$cols = @("D", "I", "J", "K", "L", "M", "O", "P")

# Snapshot the original values for every affected column/row before
# writing anything, since several rows both give and receive data.
# Value2 is used (rather than Value) to avoid the date/currency-specific
# COM variant wrapping performed by Value.
$original = @{}
foreach ($col in $cols) {
    $original[$col] = @{}
    for ($r = 2; $r -le 43; $r++) {
        $original[$col][$r] = $ws.Range("$col$r").Value2
    }
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $original[$col][$srcRow]
    }
}
